# Add new "SequencerPassthroughApp" application type entry to the
# "Application Types" sheet, right before the existing "OutputApp" entry.
#
# This inserts two rows (a header row + a description row) at row 21,
# which pushes the "OutputApp" header (previously row 21) and everything
# below it down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Application Types - Elara Appli")

# Insert two blank rows at 21/22 (old row 21 "OutputApp" header becomes row 23).
$ws.Rows("21:22").Insert()

# Pull the A:C formatting (styles) for the new rows from the rows they are
# modeled on, now shifted to rows 23 (header style) and 24 (description
# style), so the new rows look identical to every other header/description
# pair in this table.
$ws.Range("A23:C23").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)

$ws.Range("A24:C24").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)

# Match the custom row heights used by the new rows.
$ws.Rows("21:21").RowHeight = 23
$ws.Rows("22:22").RowHeight = 33

# Header row: new application type name.
$ws.Range("A21").Value = "SequencerPassthroughApp"

# Description row: "-" marker (column B) and the description text (column C),
# where the leading app-type name is italicised, matching the style used
# for every other description cell in this column.
$ws.Range("B22").Value = "-"

$descText = "SequencerProcessorApp that directly invokes a standard passthrough processor with every sequenced command;  the passthrough processor routes a single event for every command with the same payload as the command"
$ws.Range("C22").Value = $descText

$italicLen = "SequencerProcessorApp".Length
$totalLen = $descText.Length

$italicPart = $ws.Range("C22").Characters(1, $italicLen)
$italicPart.Font.Italic = $true
$italicPart.Font.Size = 11
$italicPart.Font.Name = "Helvetica Neue"
$italicPart.Font.Color = 0

$restPart = $ws.Range("C22").Characters($italicLen + 1, $totalLen - $italicLen)
$restPart.Font.Italic = $false
$restPart.Font.Size = 11
$restPart.Font.Name = "Helvetica Neue"
$restPart.Font.Color = 0
